$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- "Arrow: Right 11" shape: update sample command text
#     "add n/David …" -> "add t/Artemis …"
$arrow = $s.Shapes.Item("Right Arrow 11")
$arrow.TextFrame.TextRange.Text = "add t/Artemis …"

# --- Table 22 (AddCommand diagram, top-right in original deck)
#     toAdd = "Name: David"      -> toAdd = "Title: Artemis"
#     prevAddressBook = s2       -> prevBookShelf = s2
$tbl22 = $s.Shapes.Item("Table 22")
$trAdd = $tbl22.Table.Cell(2, 1).Shape.TextFrame.TextRange
$trAdd.Paragraphs(1).Text = "toAdd = “Title: Artemis”"
$trAdd.Paragraphs(2).Text = "prevBookShelf = s2"

# --- Table 21 (DeleteCommand diagram)
#     prevAddressBook = s3 -> prevBookShelf = s3
$tbl21 = $s.Shapes.Item("Table 21")
$trDel1 = $tbl21.Table.Cell(2, 1).Shape.TextFrame.TextRange
$trDel1.Paragraphs(2).Text = "prevBookShelf = s3"

# --- Table 23 (DeleteCommand diagram, duplicate)
#     prevAddressBook = s3 -> prevBookShelf = s3
$tbl23 = $s.Shapes.Item("Table 23")
$trDel2 = $tbl23.Table.Cell(2, 1).Shape.TextFrame.TextRange
$trDel2.Paragraphs(2).Text = "prevBookShelf = s3"
